# Copying the same change in config
# Adds a new "Complaint - creator read access" assignment rule row to
# Sheet1, mirroring the existing "Case File - creator read access" rule
# (row 26) but for the COMPLAINT object type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the formatting (styles, borders, row height, etc.) of the last
# existing rule row (26) into the new row (27) so the new row matches the
# look of the rest of the table.
$ws.Range("A26:H26").Copy($ws.Range("A27:H27"))
$ws.Rows.Item(27).RowHeight = $ws.Rows.Item(26).RowHeight

# Fill in the new rule's data:
#  B27 - rule name
#  C27 - object type the rule applies to
#  H27 - Expression 3 (participant to add)
$ws.Cells.Item(27, 2).Value = "Complaint - creator read access"
$ws.Cells.Item(27, 3).Value = "COMPLAINT"
$ws.Cells.Item(27, 8).Value = "reader, creator"

$ws.Range("A1").Select() | Out-Null
